$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.938.08"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "'1.564.44"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'206.10"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").Value = "'0.487"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'22.10"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "'0.0584"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "'1.787.05"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "'1.571.17"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'3.75"
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "'0.513"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").Value = "'26.938.38"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "'61.74"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "'214.08"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").Value = "'4.09"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'9.35"
$ws.Range("E23").Value = "  -3.10%  "
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").Value = "'151.94"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  -2.86%  "
$ws.Range("D27").Value = "'14.84"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "'1.384.87"
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("D34").Value = "'2.90"
$ws.Range("E34").Value = "  -1.29%  "
$ws.Range("E35").Value = "  +1.05%  "
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").Value = "'0.942"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E38").Value = "  -2.59%  "
$ws.Range("D39").Value = "'0.807"
$ws.Range("D40").Value = "'0.511"
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("D44").Value = "'1.79"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "'63.39"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "'1.700.02"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "'85.28"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'0.0₇0972"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'0.0946"
$ws.Range("E51").Value = "  -1.03%  "
